# Auto-generated Excel COM-interop script to apply the Zodiark_Profits.xlsx diff
# Updates numeric cell values across sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1914.2858
$ws.Range("I18").Value = 1850
$ws.Range("K18").Value = 1850
$ws.Range("M18").Value = -1566
$ws.Range("H47").Value = 50000
$ws.Range("I47").Value = 50000
$ws.Range("K47").Value = 50000
$ws.Range("M47").Value = -49028
$ws.Range("H92").Value = 1054.875
$ws.Range("J92").Value = 773
$ws.Range("L92").Value = 773
$ws.Range("N92").Value = -3269
$ws.Range("H112").Value = 8207.157999999999
$ws.Range("J112").Value = 8784.471
$ws.Range("L112").Value = 26353.413
$ws.Range("N112").Value = -28569.413
$ws.Range("H132").Value = 3509.44
$ws.Range("I132").Value = 3079.0454
$ws.Range("K132").Value = 9237.136200000001
$ws.Range("M132").Value = -6707.136200000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").Value = -120763
$ws.Range("H61").Value = 3115.0476
$ws.Range("I61").Value = 2801.1177
$ws.Range("K61").Value = 2801.1177
$ws.Range("M61").Value = -2589.1177
$ws.Range("H74").Value = 7049.8037
$ws.Range("I74").Value = 4882.2046
$ws.Range("J74").Value = 14997.667
$ws.Range("K74").Value = 4882.2046
$ws.Range("L74").Value = 14997.667
$ws.Range("M74").Value = -4008.2046
$ws.Range("N74").Value = -16745.667
$ws.Range("H77").Value = 7049.8037
$ws.Range("I77").Value = 4882.2046
$ws.Range("J77").Value = 14997.667
$ws.Range("K77").Value = 24411.023
$ws.Range("L77").Value = 74988.33499999999
$ws.Range("M77").Value = -20043.023
$ws.Range("N77").Value = -83724.33499999999
$ws.Range("H122").Value = 8240
$ws.Range("I122").Value = 8964.723
$ws.Range("J122").Value = 4978.75
$ws.Range("K122").Value = 26894.169
$ws.Range("L122").Value = 14936.25
$ws.Range("M122").Value = -24444.169
$ws.Range("N122").Value = -19836.25
$ws.Range("H136").Value = 3115.0476
$ws.Range("I136").Value = 2801.1177
$ws.Range("K136").Value = 8403.3531
$ws.Range("M136").Value = -5853.3531

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 3700
$ws.Range("I82").Value = 3700
$ws.Range("K82").Value = 3700
$ws.Range("M82").Value = -3317
$ws.Range("H85").Value = 3700
$ws.Range("I85").Value = 3700
$ws.Range("K85").Value = 3700
$ws.Range("M85").Value = -2374
$ws.Range("H105").Value = 1812.9231
$ws.Range("I105").Value = 1741
$ws.Range("J105").Value = 1974.75
$ws.Range("K105").Value = 1741
$ws.Range("L105").Value = 1974.75
$ws.Range("M105").Value = 6
$ws.Range("N105").Value = -5468.75
$ws.Range("H135").Value = 71880.59
$ws.Range("J135").Value = 71880.59
$ws.Range("L135").Value = 71880.59
$ws.Range("N135").Value = -82020.59

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 280.88235
$ws.Range("I7").Value = 86.42856999999999
$ws.Range("J7").Value = 417
$ws.Range("K7").Value = 86.42856999999999
$ws.Range("L7").Value = 417
$ws.Range("M7").Value = 26.57143000000001
$ws.Range("N7").Value = -643
$ws.Range("H22").Value = 1233.4736
$ws.Range("I22").Value = 295.86667
$ws.Range("K22").Value = 295.86667
$ws.Range("M22").Value = 54.13333
$ws.Range("H99").Value = 11856335
$ws.Range("I99").Value = 5004488
$ws.Range("K99").Value = 5004488
$ws.Range("M99").Value = -5002990
$ws.Range("H107").Value = 483.25925
$ws.Range("I107").Value = 346.83334
$ws.Range("J107").Value = 756.1111
$ws.Range("K107").Value = 346.83334
$ws.Range("L107").Value = 756.1111
$ws.Range("M107").Value = 1573.16666
$ws.Range("N107").Value = -4596.1111
$ws.Range("H122").Value = 2014.5714
$ws.Range("I122").Value = 1878.7222
$ws.Range("J122").Value = 2829.6667
$ws.Range("K122").Value = 5636.1666
$ws.Range("L122").Value = 8489.000100000001
$ws.Range("M122").Value = -3186.1666
$ws.Range("N122").Value = -13389.0001
$ws.Range("H126").Value = 11856335
$ws.Range("I126").Value = 5004488
$ws.Range("K126").Value = 15013464
$ws.Range("M126").Value = -15010994

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 1156.5555
$ws.Range("I11").Value = 581.8
$ws.Range("J11").Value = 1875
$ws.Range("K11").Value = 1745.4
$ws.Range("L11").Value = 5625
$ws.Range("M11").Value = -1605.4
$ws.Range("N11").Value = -5905
$ws.Range("H68").Value = 1836.1428
$ws.Range("J68").Value = 1977.6
$ws.Range("L68").Value = 5932.799999999999
$ws.Range("N68").Value = -7554.799999999999
$ws.Range("H71").Value = 1836.1428
$ws.Range("J71").Value = 1977.6
$ws.Range("L71").Value = 17798.4
$ws.Range("N71").Value = -25910.4
$ws.Range("H98").Value = 587.8333
$ws.Range("I98").Value = 596.6667
$ws.Range("J98").Value = 579
$ws.Range("K98").Value = 1790.0001
$ws.Range("L98").Value = 1737
$ws.Range("M98").Value = -292.0001
$ws.Range("N98").Value = -4733
$ws.Range("H107").Value = 710.0714
$ws.Range("I107").Value = 602
$ws.Range("K107").Value = 1806
$ws.Range("M107").Value = 114
$ws.Range("H122").Value = 3409.4
$ws.Range("J122").Value = 899.25
$ws.Range("L122").Value = 8093.25
$ws.Range("N122").Value = -12993.25
$ws.Range("H130").Value = 4662.25
$ws.Range("I130").Value = 4662.25
$ws.Range("K130").Value = 13986.75
$ws.Range("M130").Value = -8966.75
$ws.Range("H131").Value = 3410
$ws.Range("J131").Value = 3949.5454
$ws.Range("L131").Value = 11848.6362
$ws.Range("N131").Value = -21928.6362
$ws.Range("H132").Value = 2128.625
$ws.Range("I132").Value = 1590
$ws.Range("J132").Value = 5899
$ws.Range("K132").Value = 14310
$ws.Range("L132").Value = 53091
$ws.Range("M132").Value = -11780
$ws.Range("N132").Value = -58151
$ws.Range("H140").Value = 1579.08
$ws.Range("I140").Value = 1228.2084
$ws.Range("K140").Value = 3684.6252
$ws.Range("M140").Value = 1495.3748

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 215999.86
$ws.Range("J11").Value = 375250
$ws.Range("L11").Value = 375250
$ws.Range("N11").Value = -375528
$ws.Range("H126").Value = 5977.952
$ws.Range("I126").Value = 6203.1333
$ws.Range("J126").Value = 5415
$ws.Range("K126").Value = 18609.3999
$ws.Range("L126").Value = 16245
$ws.Range("M126").Value = -16139.3999
$ws.Range("N126").Value = -21185

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 66673960
$ws.Range("I7").Value = 100006840
$ws.Range("K7").Value = 100006840
$ws.Range("M7").Value = -100006728
$ws.Range("H45").Value = 28492
$ws.Range("I45").Value = 28492
$ws.Range("K45").Value = 28492
$ws.Range("M45").Value = -28085
$ws.Range("H126").Value = 66673960
$ws.Range("I126").Value = 100006840
$ws.Range("K126").Value = 300020520
$ws.Range("M126").Value = -300018050
$ws.Range("H136").Value = 4030.851
$ws.Range("I136").Value = 4572
$ws.Range("J136").Value = 3075.8823
$ws.Range("K136").Value = 13716
$ws.Range("L136").Value = 9227.6469
$ws.Range("M136").Value = -11166
$ws.Range("N136").Value = -14327.6469

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 52633.332
$ws.Range("I70").Value = 51450
$ws.Range("K70").Value = 51450
$ws.Range("M70").Value = -51135
$ws.Range("H73").Value = 52633.332
$ws.Range("I73").Value = 51450
$ws.Range("K73").Value = 51450
$ws.Range("M73").Value = -50358
$ws.Range("H107").Value = 1166.7742
$ws.Range("I107").Value = 1187.0625
$ws.Range("K107").Value = 3561.1875
$ws.Range("M107").Value = -1641.1875
$ws.Range("H132").Value = 3614
$ws.Range("I132").Value = 2502
$ws.Range("K132").Value = 7506
$ws.Range("M132").Value = -4976
